$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G5").Value = "已售罄"
$ws.Range("F7").Value = 2333
$ws.Range("F8").Value = 1710
$ws.Range("F11").Value = 132
$ws.Range("F14").Value = 2674
$ws.Range("F17").Value = 7142
$ws.Range("F19").Value = 7287
$ws.Range("F22").Value = 5571
$ws.Range("F23").Value = 3131
$ws.Range("F24").Value = 3502
$ws.Range("F26").Value = 244
$ws.Range("F27").Value = 196
$ws.Range("F28").Value = 1924
$ws.Range("F33").Value = 491
$ws.Range("F35").Value = 2449
$ws.Range("F36").Value = 1243
$ws.Range("F37").Value = 2801
$ws.Range("F38").Value = 46
$ws.Range("F40").Value = 172
$ws.Range("F41").Value = 401
$ws.Range("F42").Value = 1102
$ws.Range("F44").Value = 487
$ws.Range("F45").Value = 538

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 221
$ws.Range("F12").Value = 331
$ws.Range("F13").Value = 23
$ws.Range("F17").Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 83

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2333
$ws.Range("F7").Value = 1710
$ws.Range("F10").Value = 132
$ws.Range("F13").Value = 83
$ws.Range("F14").Value = 2674
$ws.Range("F16").Value = 221
$ws.Range("F19").Value = 7142
$ws.Range("F21").Value = 7287
$ws.Range("F23").Value = 5571
$ws.Range("F24").Value = 3131
$ws.Range("F25").Value = 3502
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 244
$ws.Range("F30").Value = 1924
$ws.Range("F31").Value = 20
$ws.Range("F36").Value = 491
$ws.Range("F38").Value = 2449
$ws.Range("F39").Value = 1243
$ws.Range("F41").Value = 2801
$ws.Range("F42").Value = 46
$ws.Range("F44").Value = 172
$ws.Range("F45").Value = 401
$ws.Range("F46").Value = 1102
$ws.Range("F48").Value = 487
$ws.Range("F49").Value = 538
